$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 6.663840333333333
$ws.Range("H2").Value = 19.991521
$ws.Range("I2").Value = 0.3746160267057107
$ws.Range("J2").Value = 0.3746160267057107
$ws.Range("M2").Value = 1.058456666666667
$ws.Range("N2").Value = 3.17537
$ws.Range("O2").Value = 0.00154290138396175
$ws.Range("P2").Value = 0.00154290138396175
$ws.Range("Q2").Value = 7.053386226418888
$ws.Range("R2").Value = 63.48047603777
$ws.Range("S2").Value = 0.0005779955860584928
$ws.Range("T2").Value = 0.000577995586058493

$ws.Range("G3").Value = 6.663840333333333
$ws.Range("H3").Value = 19.991521
$ws.Range("I3").Value = 0.3746160267057107
$ws.Range("J3").Value = 0.3746160267057107
$ws.Range("O3").Value = 0.000417485732185422
$ws.Range("P3").Value = 0.0004174857321854219
$ws.Range("Q3").Value = 1.908539420427445
$ws.Range("R3").Value = 17.176854783847
$ws.Range("S3").Value = 0.0001563968461976272
$ws.Range("T3").Value = 0.0001563968461976272

$ws.Range("G4").Value = 6.663840333333333
$ws.Range("H4").Value = 19.991521
$ws.Range("I4").Value = 0.3746160267057107
$ws.Range("J4").Value = 0.3746160267057107
$ws.Range("M4").Value = 395.8171083333334
$ws.Range("N4").Value = 1187.451325
$ws.Range("O4").Value = 0.5769785230476177
$ws.Range("P4").Value = 0.5769785230476177
$ws.Range("Q4").Value = 2637.662011135036
$ws.Range("R4").Value = 23738.95810021532
$ws.Range("S4").Value = 0.2161454017986278
$ws.Range("T4").Value = 0.2161454017986278

$ws.Range("G5").Value = 6.663840333333333
$ws.Range("H5").Value = 19.991521
$ws.Range("I5").Value = 0.3746160267057107
$ws.Range("J5").Value = 0.3746160267057107
$ws.Range("M5").Value = 0.173927
$ws.Range("N5").Value = 0.5217809999999999
$ws.Range("O5").Value = 0.0002535315969556132
$ws.Range("P5").Value = 0.0002535315969556133
$ws.Range("Q5").Value = 1.159021757655667
$ws.Range("R5").Value = 10.431195818901
$ws.Range("S5").Value = 0.00009497699949586549
$ws.Range("T5").Value = 0.00009497699949586551

$ws.Range("G6").Value = 6.663840333333333
$ws.Range("H6").Value = 19.991521
$ws.Range("I6").Value = 0.3746160267057107
$ws.Range("J6").Value = 0.3746160267057107
$ws.Range("M6").Value = 288.6811626666667
$ws.Range("N6").Value = 866.043488
$ws.Range("O6").Value = 0.4208075582392796
$ws.Range("P6").Value = 0.4208075582392796
$ws.Range("Q6").Value = 1923.725175251694
$ws.Range("R6").Value = 17313.52657726525
$ws.Range("S6").Value = 0.1576412554753309
$ws.Range("T6").Value = 0.1576412554753309

$ws.Range("G7").Value = 11.12461466666667
$ws.Range("H7").Value = 33.373844
$ws.Range("I7").Value = 0.6253839732942893
$ws.Range("J7").Value = 0.6253839732942893
$ws.Range("M7").Value = 1.058456666666667
$ws.Range("N7").Value = 3.17537
$ws.Range("O7").Value = 0.00154290138396175
$ws.Range("P7").Value = 0.00154290138396175
$ws.Range("Q7").Value = 11.77492255803111
$ws.Range("R7").Value = 105.97430302228
$ws.Range("S7").Value = 0.0009649057979032569
$ws.Range("T7").Value = 0.000964905797903257

$ws.Range("G8").Value = 11.12461466666667
$ws.Range("H8").Value = 33.373844
$ws.Range("I8").Value = 0.6253839732942893
$ws.Range("J8").Value = 0.6253839732942893
$ws.Range("O8").Value = 0.000417485732185422
$ws.Range("P8").Value = 0.0004174857321854219
$ws.Range("Q8").Value = 3.186115597967556
$ws.Range("R8").Value = 28.675040381708
$ws.Range("S8").Value = 0.0002610888859877948
$ws.Range("T8").Value = 0.0002610888859877947

$ws.Range("G9").Value = 11.12461466666667
$ws.Range("H9").Value = 33.373844
$ws.Range("I9").Value = 0.6253839732942893
$ws.Range("J9").Value = 0.6253839732942893
$ws.Range("M9").Value = 395.8171083333334
$ws.Range("N9").Value = 1187.451325
$ws.Range("O9").Value = 0.5769785230476177
$ws.Range("P9").Value = 0.5769785230476177
$ws.Range("Q9").Value = 4403.312808682589
$ws.Range("R9").Value = 39629.81527814329
$ws.Range("S9").Value = 0.3608331212489898
$ws.Range("T9").Value = 0.3608331212489898

$ws.Range("G10").Value = 11.12461466666667
$ws.Range("H10").Value = 33.373844
$ws.Range("I10").Value = 0.6253839732942893
$ws.Range("J10").Value = 0.6253839732942893
$ws.Range("M10").Value = 0.173927
$ws.Range("N10").Value = 0.5217809999999999
$ws.Range("O10").Value = 0.0002535315969556132
$ws.Range("P10").Value = 0.0002535315969556133
$ws.Range("Q10").Value = 1.934870855129333
$ws.Range("R10").Value = 17.413837696164
$ws.Range("S10").Value = 0.0001585545974597477
$ws.Range("T10").Value = 0.0001585545974597478

$ws.Range("G11").Value = 11.12461466666667
$ws.Range("H11").Value = 33.373844
$ws.Range("I11").Value = 0.6253839732942893
$ws.Range("J11").Value = 0.6253839732942893
$ws.Range("M11").Value = 288.6811626666667
$ws.Range("N11").Value = 866.043488
$ws.Range("O11").Value = 0.4208075582392796
$ws.Range("P11").Value = 0.4208075582392796
$ws.Range("Q11").Value = 3211.466696191986
$ws.Range("R11").Value = 28903.20026572787
$ws.Range("S11").Value = 0.2631663027639487
$ws.Range("T11").Value = 0.2631663027639487
